$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new summary row (row 49) for "Baseline 2010-18 C138", following the
# same pattern as the existing "Baseline 2010-18 C136" row (row 48).
# Duplicate row 48 (values + formatting) down into the new row 49 first,
# then overwrite the values below with the correct figures for C138.
$ws.Rows("48:48").Copy()
$ws.Rows("49:49").Insert(-4121)

$ws.Range("A49").Value = "CW3M"
$ws.Range("B49").Value = "Baseline 2010-18 C138"
$ws.Range("C49").Value = "2010-18"
$ws.Range("D49").Value = 1187.0067003333331
$ws.Range("E49").Value = 1901.5157334444443
$ws.Range("F49").Value = 0.97970299999999988
$ws.Range("G49").Value = 280.33542888888883
$ws.Range("H49").Value = 9.775355222222224
$ws.Range("I49").Value = 5.3870271111111121
$ws.Range("J49").Value = 8.145128999999999
$ws.Range("K49").Value = 645.93713388888887
$ws.Range("L49").Value = 83.47062044444445
$ws.Range("M49").Value = 1455.6553682222225
$ws.Range("N49").Value = 1191.1764458888888
$ws.Range("O49").Value = 4661.9885253333332
$ws.Range("P49").Value = 27227.338324888889
$ws.Range("Q49").Value = -0.61525011111111105
$ws.Range("R49").Value = -0.00017433333333333333
$ws.Range("S49").Value = "2010-18"

# The author's last selection, recorded in the saved workbook view state.
$ws.Range("F54").Select()
